# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New computed values for column G ("K") for rows 2-12
$newValues = @{
    2  = 2
    3  = 1
    4  = 2
    5  = 1
    6  = 2
    7  = 4
    8  = 2
    9  = 1
    10 = 3
    11 = 3
    12 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
